# rnaSample_1874.xlsx — "fixed harvester column in rnasamples -- holly added
# S.GISH to harvester in bioSamples"
#
# The harvester column (column B) incorrectly held the same value as the
# rnaPreparer column ("Retrofitted_1874") for every data row. Holly's fix
# replaces the harvester entry for every sample row with "S.GISH".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2-43 (row 1 is the header).
$ws.Range("B2:B43").Value = "S.GISH"

# Column B was widened slightly (matches the new, slightly longer header
# selection width Excel remembers after the edit).
$ws.Columns("B").ColumnWidth = 7.996666666666667

# Leave the column selected, mirroring the state Excel saved after the
# harvester values were entered down column B.
$ws.Columns("B").Select()
